# Apply the "dated warning" / reviewer update to the capacity-planning workbook.
#
# Semantic changes (per the target diff):
#  1. "Server Scenarios"!C6 (Future_Users / "Maximum expected users")  2000 -> 0
#  2. "Server Scenarios"!C7 (Current_Users / "Current users")          200  -> 0
#  3. "Server Scenarios"!B2  "Version 2012.02.24" -> "Version 2012.02.24 (DATED!)"
#  4. "Server Scenarios"!D9  "Review the Real-World "Beef" Factor section in the
#     TFS Project Guidance document." -> "Safety growth factor to add to
#     proposed arcjitecture."
#  5. Active sheet / selection: "Server Scenarios" becomes the selected tab
#     (was "Hardware Configurations"); selection on "Server Scenarios" moves
#     to Q9; selection on "Hardware Configurations" moves to B2:M2.
#
# All of the other cell-value churn visible in the raw OOXML diff (I13, C14,
# C15, D16, E16, C18, I17:I21, J8:M10, D15:E17, rows 23-32 on the Hardware
# Configurations sheet, etc.) is pure formula recalculation fallout from
# changes 1 & 2 above, so it is left for Excel's own calc engine.

$wb = $excel.ActiveWorkbook

$wsScenarios = $wb.Worksheets.Item("Server Scenarios")
$wsHardware  = $wb.Worksheets.Item("Hardware Configurations")

# 1 & 2: zero out current/future user counts.
$wsScenarios.Range("C6").Value = 0
$wsScenarios.Range("C7").Value = 0

# 3: mark the version text as dated.
$wsScenarios.Range("B2").Value = "Version 2012.02.24 (DATED!)"

# 4: replace the stale "review the guidance doc" note with the new
#    safety-growth-factor explanation next to the Beef Factor input.
$wsScenarios.Range("D9").Value = "Safety growth factor to add to proposed arcjitecture."

# Recalculate so dependent formulas (Hardware Configurations' pull-through of
# B2, and all the Yes/No + capacity formulas on both sheets) pick up the new
# inputs before the selection/activation below.
$excel.Calculate()

# 5: make "Server Scenarios" the active/selected sheet again, with the
#    selections observed in the target workbook.
$wsHardware.Range("B2:M2").Select() | Out-Null
$wsScenarios.Activate()
$wsScenarios.Range("Q9").Select() | Out-Null
